# Weekly update: insert a new daily price record for "Choclo" at
# Terminal Hortofrutícola Agro Chillán, pushing the existing rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right above the current first data row of this block (row 279),
# shifting every following row down by one (279->280, ..., 357->358).
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row with the new week's record.
$ws.Range("A279").Value = 7
$ws.Range("B279").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C279").Value = 'Ñuble'
$ws.Range("D279").Value = 45204
$ws.Range("E279").Value = 16
$ws.Range("F279").Value = 100112024
$ws.Range("G279").Value = 'Choclo'
$ws.Range("H279").Value = 'Dulce o Americano'
$ws.Range("I279").Value = 'Primera'
$ws.Range("J279").Value = 40
$ws.Range("K279").Value = 25000
$ws.Range("L279").Value = 25000
$ws.Range("M279").Value = 25000
$ws.Range("N279").Value = '$/malla 70 unidades'
$ws.Range("O279").Value = 'Región de Arica y Parinacota'
$ws.Range("P279").Value = 357
$ws.Range("Q279").Value = 70
$ws.Range("R279").Value = 'Hortaliza'
